$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (changed) date column C was bumped from 2023-10-09 (45208)
# to 2023-10-13 (45212) for every data row (2..62).
$ws.Range("C2:C62").Value2 = 45212

# Row 2's file-link hyperlink formulas (columns S, T, V, W, X, Y) gained
# more specific filenames (and Y2 also picked up a "ti,llsynsmail" typo
# in its folder name, matching the upstream change exactly).
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/artfynd/A 31987-2023 artfynd.xlsx", "A 31987-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/kartor/A 31987-2023 karta.png", "A 31987-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/klagomål/A 31987-2023 fsc-klagomål.docx", "A 31987-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/klagomålsmail/A 31987-2023 fsc-klagomål mail.docx", "A 31987-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/tillsyn/A 31987-2023 tillsynsbegäran.docx", "A 31987-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/ti,llsynsmail/A 31987-2023 tillsynsbegäran mail.docx", "A 31987-2023")'
